# Regenerate save_data: write new "K" (Strike#) values for rows 2-57 in column G.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values (column G), one per row from row 2 to row 57 (56 values).
$newK = @(2,1,1,2,2,0,0,0,0,0,3,2,1,0,1,0,2,1,1,0,2,0,0,0,1,1,1,2,1,2,1,0,1,1,3,1,0,1,2,1,1,1,2,1,0,0,0,1,2,0,2,1,3,1,1,0)

$startRow = 2
for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}

$wb.Save()
